# petty-cashBook-2021.xlsx — "Update 30-Mar-2021, midday update."
# Applies the recorded edits to Sheet1 of the active workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# -------------------------------------------------------------------------
# Opening balance (row 2)
# -------------------------------------------------------------------------
$ws.Range("E2").Value = 838525

# -------------------------------------------------------------------------
# Row 3 — date 44284 (30-Mar-2021), Wages Expense debit recomputed
# -------------------------------------------------------------------------
$ws.Range("A3").Value = 44284
$ws.Range("D3").Formula = "=60000+260000"

# -------------------------------------------------------------------------
# Row 4 — now "BELI kresek"
# -------------------------------------------------------------------------
$ws.Range("B4").Value = "BELI kresek"
$ws.Range("D4").Formula = "=54000"

# -------------------------------------------------------------------------
# Row 5 — now "TRANSFER BCA" (debit, was credit "A/R")
# -------------------------------------------------------------------------
$ws.Range("B5").Value = "TRANSFER BCA"
$ws.Range("C5").Clear()
$ws.Range("D5").Formula = "=1207500+19680000+577500+368500+400000"

# -------------------------------------------------------------------------
# Row 6 — now "A/R" (credit, was debit "TRANSFER DANAMON")
# -------------------------------------------------------------------------
$ws.Range("B6").Value = "A/R"
$ws.Range("D6").Clear()
$ws.Range("C6").Formula = "=45000000+6000000+29663500"

# -------------------------------------------------------------------------
# Row 7 — now "TRANSFER DANAMON" (debit formula, was literal "A/P" amount)
# -------------------------------------------------------------------------
$ws.Range("B7").Value = "TRANSFER DANAMON"
$ws.Range("D7").Formula = "=45000000"

# -------------------------------------------------------------------------
# Row 8 — "SALES - cash/retail" unchanged, credit formula updated
# -------------------------------------------------------------------------
$ws.Range("C8").Formula = "=21643525+15508975-29663500"

# -------------------------------------------------------------------------
# Row 9 — label now "SELISIH - kurang", debit literal updated
# -------------------------------------------------------------------------
$ws.Range("B9").Value = "SELISIH - kurang"
$ws.Range("D9").Value = 135000

# -------------------------------------------------------------------------
# Row 10 — label now "SETOR KE BANK" (debit formula, was credit literal)
# -------------------------------------------------------------------------
$ws.Range("B10").Value = "SETOR KE BANK"
$ws.Range("C10").Clear()
$ws.Range("D10").Formula = "=21000000"

# -------------------------------------------------------------------------
# Row 11 — new date 44285 (31-Mar-2021) starts here, "Wages Expense"
# -------------------------------------------------------------------------
$ws.Range("A11").Value = 44285
$ws.Range("B11").Value = "Wages Expense"
$ws.Range("D11").Formula = "=60000"

# -------------------------------------------------------------------------
# Row 12 — date moved out (now blank), label now "A/R" (credit formula)
# -------------------------------------------------------------------------
$ws.Range("A12").Clear()
$ws.Range("B12").Value = "A/R"
$ws.Range("D12").Clear()
$ws.Range("C12").Formula = "=16250000"

# -------------------------------------------------------------------------
# Row 13 — label now "TRANSFER BCA" (debit formula, was credit "A/R")
# -------------------------------------------------------------------------
$ws.Range("B13").Value = "TRANSFER BCA"
$ws.Range("C13").Clear()
$ws.Range("D13").Formula = "=16250000+910000+1458000"

# -------------------------------------------------------------------------
# Rows 14-36 — all the old entries for these dates are gone; only the
# running E-column balance formulas remain (they recalc automatically).
# -------------------------------------------------------------------------
$ws.Range("B14").Clear()
$ws.Range("D14").Clear()

$ws.Range("B15").Clear()
$ws.Range("D15").Clear()

$ws.Range("B16").Clear()
$ws.Range("D16").Clear()

$ws.Range("B17").Clear()
$ws.Range("C17").Clear()

$ws.Range("B18").Clear()
$ws.Range("D18").Clear()

$ws.Range("A19").Clear()
$ws.Range("B19").Clear()
$ws.Range("D19").Clear()

$ws.Range("B20").Clear()
$ws.Range("D20").Clear()

$ws.Range("B21").Clear()
$ws.Range("D21").Clear()

$ws.Range("B22").Clear()
$ws.Range("C22").Clear()

$ws.Range("B23").Clear()
$ws.Range("C23").Clear()

$ws.Range("B24").Clear()
$ws.Range("D24").Clear()

$ws.Range("B25").Clear()
$ws.Range("D25").Clear()

$ws.Range("A26").Clear()
$ws.Range("B26").Clear()
$ws.Range("D26").Clear()

$ws.Range("B27").Clear()
$ws.Range("D27").Clear()

$ws.Range("B28").Clear()
$ws.Range("D28").Clear()

$ws.Range("B29").Clear()
$ws.Range("C29").Clear()

$ws.Range("B30").Clear()
$ws.Range("C30").Clear()

$ws.Range("B31").Clear()
$ws.Range("D31").Clear()

$ws.Range("A32").Clear()
$ws.Range("B32").Clear()
$ws.Range("D32").Clear()

$ws.Range("B33").Clear()
$ws.Range("D33").Clear()

$ws.Range("B34").Clear()
$ws.Range("D34").Clear()

$ws.Range("B35").Clear()
$ws.Range("D35").Clear()

$ws.Range("B36").Clear()
$ws.Range("D36").Clear()

# -------------------------------------------------------------------------
# View: scrolled back to the top, active cell now D14
# -------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$ws.Range("D14").Select()
